$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three validation strings (VT187-1368, VT187-1369, VT187-1370) so that
# validate_AppMinimized points to "com.symbol.enterprisebrowser" instead of "AppsScreen".
# These live in cells H4 (1368), H5 (1369) and H9 (1370) of the TestCases sheet.

$text1368 = "validate1`n{`nvalidate_PageTitle=RE 2.2 Tests`n};`nvalidate2`n{`nvalidate_PageTitle=Device & Application`n};`nvalidate3`n{`nvalidate_Text_Exists=VT187-1368`n};`nvalidate4`n{`nvalidate_AppMinimized=com.symbol.enterprisebrowser`n};"

$text1369 = "validate1`n{`nvalidate_PageTitle=RE 2.2 Tests`n};`nvalidate2`n{`nvalidate_PageTitle=Device & Application`n};`nvalidate3`n{`nvalidate_Text_Exists=VT187-1369`n};`nvalidate4`n{`nvalidate_AppMinimized=com.symbol.enterprisebrowser`n};`nvalidate5`n{`nvalidate_Result=Minimized`nvalidate_Result=Restored`n};"

$text1370 = "validate1`n{`nvalidate_PageTitle=RE 2.2 Tests`n};`nvalidate2`n{`nvalidate_PageTitle=Device & Application`n};`nvalidate3`n{`nvalidate_Text_Exists=VT187-1370`n};`nvalidate4`n{`nvalidate_AppMinimized=com.symbol.enterprisebrowser`n};`nvalidate5`n{`nvalidate_PageTitle=Device & Application`nvalidate_Result=Minimized`nvalidate_Result=Restored`n};`n"

$ws.Range("H4").Value = $text1368
$ws.Range("H5").Value = $text1369
$ws.Range("H9").Value = $text1370

# Those rows grew taller once the AppMinimized text became longer.
$ws.Rows.Item(4).RowHeight = 217.5
$ws.Rows.Item(5).RowHeight = 281.25
$ws.Rows.Item(9).RowHeight = 306.75

# Restore the selected cell shown when the workbook was last saved.
$null = $ws.Range("G2").Select()
